# Reduce exam duration from 3 hours (180 min) to 2 hours (120 min) and
# shorten every time slot's end time by one hour accordingly.

$wb = $excel.ActiveWorkbook

function Convert-TimeSlot($slot) {
    switch ($slot) {
        "09:00 - 12:00" { return "09:00 - 11:00" }
        "14:00 - 17:00" { return "14:00 - 16:00" }
        default { return $slot }
    }
}

# --- Exam_Schedule sheet -------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Exam_Schedule")
$lastRowSchedule = $wsSchedule.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRowSchedule; $row++) {
    # duration (column D)
    $wsSchedule.Cells.Item($row, 4).Value2 = "2 hours"
    # duration_minutes (column E)
    $wsSchedule.Cells.Item($row, 5).Value2 = 120
    # time_slot (column K)
    $oldSlot = $wsSchedule.Cells.Item($row, 11).Value2
    $wsSchedule.Cells.Item($row, 11).Value2 = Convert-TimeSlot $oldSlot
}

# --- Exam_Classrooms sheet ------------------------------------------------
$wsClassrooms = $wb.Worksheets.Item("Exam_Classrooms")
$lastRowClassrooms = $wsClassrooms.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRowClassrooms; $row++) {
    # Time Slot (column E)
    $oldSlot = $wsClassrooms.Cells.Item($row, 5).Value2
    $wsClassrooms.Cells.Item($row, 5).Value2 = Convert-TimeSlot $oldSlot
    # Duration (column I)
    $wsClassrooms.Cells.Item($row, 9).Value2 = "2 hours"
}

# --- Configuration sheet ---------------------------------------------------
$wsConfig = $wb.Worksheets.Item("Configuration")
$wsConfig.Cells.Item(3, 2).Value2 = 120

# --- Department_Summary sheet ----------------------------------------------
$wsDept = $wb.Worksheets.Item("Department_Summary")
$lastRowDept = $wsDept.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRowDept; $row++) {
    $numExams = $wsDept.Cells.Item($row, 2).Value2
    $newTotalMinutes = $numExams * 120
    $wsDept.Cells.Item($row, 3).Value2 = $newTotalMinutes
    $wsDept.Cells.Item($row, 5).Value2 = $newTotalMinutes / 60
}

Write-Host "Applied B4 elective duration change (3 hours -> 2 hours)"
